$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Do you want to show the operator name on a buyer research ?" column
# (column F) is removed entirely: delete the whole column so everything to
# its right (columns G onward) shifts one place to the left, the sheet's
# used range shrinks by one column, and the now-orphaned shared string is
# dropped automatically on save.
$ws.Columns("F").Delete() | Out-Null

# After deleting column F, Excel leaves the active cell at the position
# where the deleted column used to be.
$ws.Range("F4").Select() | Out-Null
